# Update the cryptocurrency price (D) and 1h volume/change (E) columns
# with freshly scraped values, matching the upstream GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.600.28"
$ws.Range("E2").Value = "  +4.19%  "

$ws.Range("D3").Value = "3.149.86"
$ws.Range("E3").Value = "  +2.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.96"
$ws.Range("E5").Value = "  +4.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "624.51"
$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.267"
$ws.Range("E7").Value = "  +24.32%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  +5.79%  "

$ws.Range("D10").Value = "3.144.73"
$ws.Range("E10").Value = "  +2.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.581"
$ws.Range("E11").Value = "  +30.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000251"
$ws.Range("E12").Value = "  +27.45%  "

$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("D14").Value = "3.720.18"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.22"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.33"
$ws.Range("E16").Value = "  +6.88%  "

$ws.Range("D17").Value = "79.384.85"
$ws.Range("E17").Value = "  +4.08%  "

$ws.Range("D18").Value = "3.135.53"
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.17"
$ws.Range("E19").Value = "  +4.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.97"
$ws.Range("E20").Value = "  +15.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.48"
$ws.Range("E21").Value = "  +12.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.03"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.16"
$ws.Range("E23").Value = "  +14.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.88"
$ws.Range("E24").Value = "  +6.63%  "

$ws.Range("D25").Value = "3.300.42"
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "75.61"
$ws.Range("E26").Value = "  +4.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.65"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.73"
$ws.Range("E28").Value = "  +6.49%  "

$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000121"
$ws.Range("E30").Value = "  +11.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.87"
$ws.Range("E32").Value = "  +6.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "550.22"
$ws.Range("E33").Value = "  +9.78%  "

$ws.Range("E34").Value = "  +2.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.149"
$ws.Range("E35").Value = "  +15.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.98"
$ws.Range("E36").Value = "  +2.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.89"
$ws.Range("E37").Value = "  +9.61%  "

$ws.Range("E38").Value = "  +18.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.400"
$ws.Range("E40").Value = "  +5.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.82"
$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.55"
$ws.Range("E44").Value = "  +7.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "187.30"
$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.80"
$ws.Range("E46").Value = "  +7.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  +9.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.777"
$ws.Range("E48").Value = "  -2.51%  "

$ws.Range("E49").Value = "  +1.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.42"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.19"
$ws.Range("E51").Value = "  +6.37%  "
